$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 427, shifting rows 427:510 down to 428:511
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with its data
$ws.Cells.Item(427, 1).Value = 4
$ws.Cells.Item(427, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(427, 3).Value = "Los Lagos"
$ws.Cells.Item(427, 4).Value = 45275
$ws.Cells.Item(427, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(427, 5).Value = 10
$ws.Cells.Item(427, 6).Value = "Fruta"
$ws.Cells.Item(427, 7).Value = 100108
$ws.Cells.Item(427, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(427, 9).Value = 100108005
$ws.Cells.Item(427, 10).Value = "Piña"
$ws.Cells.Item(427, 11).Value = "Caramelo"
$ws.Cells.Item(427, 12).Value = "Primera"
$ws.Cells.Item(427, 13).Value = 100
$ws.Cells.Item(427, 14).Value = 26000
$ws.Cells.Item(427, 15).Value = 26000
$ws.Cells.Item(427, 16).Value = 26000
$ws.Cells.Item(427, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(427, 18).Value = "Ecuador"
$ws.Cells.Item(427, 19).Value = 1857
$ws.Cells.Item(427, 20).Value = 14
